# Auto-update draw results: append the newest Pick 3 draw as a new row
# at the bottom of the "Results" sheet (mirrors the nightly scraper job
# that produces this workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data block.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New draw result row.
$date        = "2025-12-05"
$game        = "Pick 3"
$phase       = "251205"
$result      = "8-9-0"
$insertedAt  = "2025-12-05T21:40:29.694+04:00"

$rng = $ws.Range("A" + $newRow + ":E" + $newRow)

# Every column in this sheet stores plain text (dates, phase codes and
# timestamps are kept as literal strings, not Excel dates/numbers), so
# force text entry with a leading apostrophe for values that would
# otherwise be auto-converted to a date serial / number.
$ws.Range("A" + $newRow).Value = "'" + $date
$ws.Range("B" + $newRow).Value = $game
$ws.Range("C" + $newRow).Value = "'" + $phase
$ws.Range("D" + $newRow).Value = $result
$ws.Range("E" + $newRow).Value = $insertedAt

# Drop the quote-prefix formatting the apostrophe trick introduces so the
# new row keeps the same (default/unstyled) look as every other row.
$rng.Style = "Normal"
